# Update the LR-pairs data table on the active sheet to reflect the
# revised NATMI output (per Dr Hou's advice): recomputed statistics for the
# existing Sending-cluster/Target-cluster combinations, plus two new target
# clusters ("M1" and "Neutro") for every sending cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D hold text labels (Sending cluster, Ligand symbol, Receptor
# symbol, Target cluster); columns E-T hold the recomputed numeric stats.
$data = @(
    @("ECs", "Bmp4", "Rgmb", "ECs", 3, 1, 5.286854333333333, 15.860563, 0.4075099150374557, 0.4075099150374557, 3, 1, 8.646358666666666, 25.939076, 0.1920122657718527, 0.1920122657718527, 45.71203878442088, 411.4083490597879, 0.07824690211083705, 0.07824690211083704),
    @("ECs", "Bmp4", "Rgmb", "FAPs", 3, 1, 5.286854333333333, 15.860563, 0.4075099150374557, 0.4075099150374557, 3, 1, 7.131051333333333, 21.393154, 0.1583613838652608, 0.1583613838652608, 37.70082964285577, 339.307466785702, 0.06453383408414633, 0.06453383408414633),
    @("ECs", "Bmp4", "Rgmb", "M1", 3, 1, 5.286854333333333, 15.860563, 0.4075099150374557, 0.4075099150374557, 2, 0.6666666666666666, 1.531490333333333, 4.594471, 0.0340102626143302, 0.03401026261433021, 8.096766305241445, 72.87089674717301, 0.01385951922836726, 0.01385951922836726),
    @("ECs", "Bmp4", "Rgmb", "M2", 3, 1, 5.286854333333333, 15.860563, 0.4075099150374557, 0.4075099150374557, 3, 1, 4.621346333333332, 13.864039, 0.1026276163861554, 0.1026276163861554, 24.43238488821744, 219.891463993957, 0.04182177123401878, 0.04182177123401878),
    @("ECs", "Bmp4", "Rgmb", "Neutro", 3, 1, 5.286854333333333, 15.860563, 0.4075099150374557, 0.4075099150374557, 2, 0.6666666666666666, 0.6145123333333333, 1.843537, 0.0136466586706575, 0.01364665867065751, 3.248837192370111, 29.239534731331, 0.005561148715424798, 0.005561148715424798),
    @("ECs", "Bmp4", "Rgmb", "sCs", 3, 1, 5.286854333333333, 15.860563, 0.4075099150374557, 0.4075099150374557, 3, 1, 22.48548233333333, 67.456447, 0.4993418126917433, 0.4993418126917434, 118.8774697110734, 1069.897227399661, 0.2034867396646614, 0.2034867396646614),
    @("FAPs", "Bmp4", "Rgmb", "ECs", 3, 1, 3.035834666666667, 9.107504, 0.2340016669801247, 0.2340016669801247, 3, 1, 8.646358666666666, 25.939076, 0.1920122657718527, 0.1920122657718527, 26.24891538070044, 236.240238426304, 0.04493119027124427, 0.04493119027124427),
    @("FAPs", "Bmp4", "Rgmb", "FAPs", 3, 1, 3.035834666666667, 9.107504, 0.2340016669801247, 0.2340016669801247, 3, 1, 7.131051333333333, 21.393154, 0.1583613838652608, 0.1583613838652608, 21.64869284751289, 194.838235627616, 0.03705682780975045, 0.03705682780975045),
    @("FAPs", "Bmp4", "Rgmb", "M1", 3, 1, 3.035834666666667, 9.107504, 0.2340016669801247, 0.2340016669801247, 2, 0.6666666666666666, 1.531490333333333, 4.594471, 0.0340102626143302, 0.03401026261433021, 4.649351445598223, 41.84416301038401, 0.007958458146185083, 0.007958458146185085),
    @("FAPs", "Bmp4", "Rgmb", "M2", 3, 1, 3.035834666666667, 9.107504, 0.2340016669801247, 0.2340016669801247, 3, 1, 4.621346333333332, 13.864039, 0.1026276163861554, 0.1026276163861554, 14.02964340540622, 126.266790648656, 0.02401503331255713, 0.02401503331255713),
    @("FAPs", "Bmp4", "Rgmb", "Neutro", 3, 1, 3.035834666666667, 9.107504, 0.2340016669801247, 0.2340016669801247, 2, 0.6666666666666666, 0.6145123333333333, 1.843537, 0.0136466586706575, 0.01364665867065751, 1.865557844627556, 16.790020601648, 0.003193340877642629, 0.003193340877642629),
    @("FAPs", "Bmp4", "Rgmb", "sCs", 3, 1, 3.035834666666667, 9.107504, 0.2340016669801247, 0.2340016669801247, 3, 1, 22.48548233333333, 67.456447, 0.4993418126917433, 0.4993418126917434, 68.26220676425422, 614.359860878288, 0.1168468165627451, 0.1168468165627451),
    @("sCs", "Bmp4", "Rgmb", "ECs", 3, 1, 4.650871, 13.952613, 0.3584884179824196, 0.3584884179824196, 3, 1, 8.646358666666666, 25.939076, 0.1920122657718527, 0.1920122657718527, 40.21309877839866, 361.9178890055879, 0.06883417338977137, 0.06883417338977137),
    @("sCs", "Bmp4", "Rgmb", "FAPs", 3, 1, 4.650871, 13.952613, 0.3584884179824196, 0.3584884179824196, 3, 1, 7.131051333333333, 21.393154, 0.1583613838652608, 0.1583613838652608, 33.16559984571133, 298.490398611402, 0.056770721971364, 0.05677072197136402),
    @("sCs", "Bmp4", "Rgmb", "M1", 3, 1, 4.650871, 13.952613, 0.3584884179824196, 0.3584884179824196, 2, 0.6666666666666666, 1.531490333333333, 4.594471, 0.0340102626143302, 0.03401026261433021, 7.122763978080333, 64.104875802723, 0.01219228523977786, 0.01219228523977787),
    @("sCs", "Bmp4", "Rgmb", "M2", 3, 1, 4.650871, 13.952613, 0.3584884179824196, 0.3584884179824196, 3, 1, 4.621346333333332, 13.864039, 0.1026276163861554, 0.1026276163861554, 21.49328564265633, 193.439570783907, 0.0367908118395795, 0.0367908118395795),
    @("sCs", "Bmp4", "Rgmb", "Neutro", 3, 1, 4.650871, 13.952613, 0.3584884179824196, 0.3584884179824196, 2, 0.6666666666666666, 0.6145123333333333, 1.843537, 0.0136466586706575, 0.01364665867065751, 2.858017590242333, 25.722158312181, 0.004892169077590078, 0.004892169077590079),
    @("sCs", "Bmp4", "Rgmb", "sCs", 3, 1, 4.650871, 13.952613, 0.3584884179824196, 0.3584884179824196, 3, 1, 22.48548233333333, 67.456447, 0.4993418126917433, 0.4993418126917434, 104.5770777051123, 941.1936993460109, 0.1790082564643368, 0.1790082564643368)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
